# "adds Cost to results dataset"
# Insert a new data row (row 4) in the "results" sheet for the "Cost" metric,
# right after the existing "LCOD" row, pushing the remaining metric rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")
$ws.Activate()

# Insert a new row at position 4, shifting rows 4-8 down to 5-9.
$ws.Rows.Item(4).Insert()

# Fill in the new row: same Technology as the other rows, Variable/Index = "Cost".
$ws.Cells.Item(4, 1).Value = $ws.Cells.Item(5, 1).Value2
$ws.Cells.Item(4, 2).Value = "Cost"
$ws.Cells.Item(4, 3).Value = "Cost"

# Match the author's resulting selection/active cell.
$ws.Range("D4").Select()

# Resize column A to fit its (now slightly different) contents.
$ws.Columns.Item(1).EntireColumn.AutoFit()
